$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Replace the italic closing paragraph's text with the new image
#    prompt text. Do this first, before any other text is inserted
#    that might coincidentally contain the same phrase.
# ------------------------------------------------------------------
$old = "Check out our review of Cleopatra slot - a classic 5-reel, 20-payline slot game by IGT. Play for free or with real money and win big with bonuses!"
$new = "Create a feature image for Cleopatra that fits the theme of the game. The image should be in a cartoon style and should showcase a happy Maya warrior with glasses. The image should be eye-catching and draw the attention of potential players to the game. The Maya warrior should be wearing traditional clothing and accessories and should be holding a symbol from the game, such as a Scarab or a Cartouche. The background should feature elements of ancient Egyptian culture, such as hieroglyphics or pyramids, to tie in with the theme of the game. The colors used should be bright and bold to make the image stand out."

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# ------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Cleopatra Slot for Free -
#    Review and Bonuses" paragraph near the end of the document
#    (leave the real Heading1 title paragraph alone).
# ------------------------------------------------------------------
$target = "Play Cleopatra Slot for Free - Review and Bonuses"
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
  $p = $d.Paragraphs.Item($i)
  $text = $p.Range.Text.TrimEnd([char]13)
  $styleName = $p.Range.ParagraphFormat.Style.NameLocal
  if ($text -eq $target -and $styleName -ne "Heading 1") {
    $p.Range.Delete()
  }
}

# ------------------------------------------------------------------
# 3) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)            # wdCollapseEnd
$titleRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(2)
$newPara.Style = "Normal"          # drop the inherited Heading1 style

$newRange = $newPara.Range
$newRange.Collapse(1)              # wdCollapseStart
$startPos = $newRange.Start
$newRange.InsertAfter("Meta description: Check out our review of Cleopatra slot - a classic 5-reel, 20-payline slot game by IGT. Play for free or with real money and win big with bonuses!")

$boldLen = ("Meta description").Length
$boldRange = $d.Range($startPos, $startPos + $boldLen)
$boldRange.Font.Bold = 1

Write-Output "Done"
